$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.681.69"
$ws.Range("E2").Value = "'  -4.16%  "
$ws.Range("D3").Value = "'2.278.62"
$ws.Range("E3").Value = "'  -5.35%  "
$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("D5").Value = "'538.65"
$ws.Range("E5").Value = "'  -3.75%  "
$ws.Range("D6").Value = "'130.42"
$ws.Range("E6").Value = "'  -3.71%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "'  +0.09%  "
$ws.Range("D8").Value = "'0.566"
$ws.Range("E8").Value = "'  -3.57%  "
$ws.Range("D9").Value = "'2.276.02"
$ws.Range("E9").Value = "'  -5.40%  "
$ws.Range("E10").Value = "'  -5.49%  "
$ws.Range("E11").Value = "'  -3.33%  "
$ws.Range("E12").Value = "'  +0.28%  "
$ws.Range("E13").Value = "'  -5.71%  "
$ws.Range("D14").Value = "'23.37"
$ws.Range("E14").Value = "'  -5.52%  "
$ws.Range("D15").Value = "'2.688.40"
$ws.Range("E15").Value = "'  -5.31%  "
$ws.Range("D16").Value = "'57.709.03"
$ws.Range("E16").Value = "'  -3.98%  "
$ws.Range("E17").Value = "'  -4.76%  "
$ws.Range("D18").Value = "'2.287.18"
$ws.Range("E18").Value = "'  -5.69%  "
$ws.Range("D19").Value = "'10.54"
$ws.Range("E19").Value = "'  -6.16%  "
$ws.Range("D20").Value = "'4.23"
$ws.Range("E20").Value = "'  -6.55%  "
$ws.Range("D21").Value = "'311.72"
$ws.Range("E21").Value = "'  -4.45%  "
$ws.Range("E22").Value = "'  -6.02%  "
$ws.Range("E23").Value = "'  -0.10%  "
$ws.Range("D24").Value = "'62.71"
$ws.Range("E24").Value = "'  -3.10%  "
$ws.Range("D25").Value = "'0.165"
$ws.Range("E25").Value = "'  -3.34%  "
$ws.Range("E26").Value = "'  +0.20%  "
$ws.Range("D27").Value = "'7.93"
$ws.Range("E27").Value = "'  -7.16%  "
$ws.Range("E28").Value = "'  -7.90%  "
$ws.Range("E29").Value = "'  -4.73%  "
$ws.Range("D30").Value = "'169.64"
$ws.Range("E30").Value = "'  -0.67%  "
$ws.Range("D31").Value = "'0.0₃0716"
$ws.Range("E31").Value = "'  -6.92%  "
$ws.Range("D32").Value = "'1.08"
$ws.Range("E32").Value = "'  +0.28%  "
$ws.Range("D33").Value = "'5.70"
$ws.Range("E33").Value = "'  -6.88%  "
$ws.Range("E34").Value = "'  -5.98%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("D36").Value = "'17.65"
$ws.Range("E37").Value = "'  +0.05%  "
$ws.Range("E38").Value = "'  -8.23%  "
$ws.Range("D39").Value = "'3.88"
$ws.Range("E39").Value = "'  -7.12%  "
$ws.Range("D40").Value = "'37.83"
$ws.Range("E40").Value = "'  -1.86%  "
$ws.Range("E41").Value = "'  -7.54%  "
$ws.Range("B42").Value = "'Aave"
$ws.Range("C42").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'139.45"
$ws.Range("E42").Value = "'  -6.27%  "
$ws.Range("B43").Value = "'Bittensor"
$ws.Range("C43").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "'286.62"
$ws.Range("E43").Value = "'  -11.83%  "
$ws.Range("E44").Value = "'  -5.26%  "
$ws.Range("D45").Value = "'0.0946"
$ws.Range("E45").Value = "'  -2.45%  "
$ws.Range("E46").Value = "'  -3.68%  "
$ws.Range("D47").Value = "'0.548"
$ws.Range("E47").Value = "'  -4.85%  "
$ws.Range("D48").Value = "'18.11"
$ws.Range("E48").Value = "'  -9.30%  "
$ws.Range("D49").Value = "'0.0210"
$ws.Range("E49").Value = "'  -5.09%  "
$ws.Range("E50").Value = "'  -0.89%  "
$ws.Range("D51").Value = "'16.42"
$ws.Range("E51").Value = "'  -4.12%  "
